$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "1 (set)"/"1 (set)"/"2 (sets)" text values replaced with plain numbers.
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 3

# New row 44: "Isolation Feet" part entry.
$ws.Range("A44").Value = "Isolation Feet"
$ws.Range("C44").Value = 4
$ws.Range("D44").Value = 4
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 14.2
$ws.Range("G44").Value = "Amazon"

# Match formatting of the row above (style carries over font/number format/etc.)
$ws.Range("A43").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("C43:G43").Copy()
$ws.Range("C44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the column-specific number formats (style carried over from A43/G43
# is for column A / column G; fix numeric columns back to their own styles).
$ws.Range("C44:E44").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C44").Value = 4
$ws.Range("D44").Value = 4
$ws.Range("E44").Value = 0

$ws.Range("C2:E2").Copy()
$ws.Range("C44").PasteSpecial(-4122)
$ws.Range("C44").Value = 4
$ws.Range("D44").Value = 4
$ws.Range("E44").Value = 0

$ws.Range("F2").Copy()
$ws.Range("F44").PasteSpecial(-4122)
$ws.Range("F44").Value = 14.2

$ws.Range("G43").Copy()
$ws.Range("G44").PasteSpecial(-4122)
$ws.Range("G44").Value = "Amazon"
$excel.CutCopyMode = $false

# New hyperlink for the added part.
$ws.Hyperlinks.Add($ws.Range("A44"), "https://www.amazon.com/s?k=isolation+feet")

# Restore column-A cell style (the hyperlink add can stamp a fresh style index).
$ws.Range("A43").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the sheet's dimension/view to reflect the new last row.
$ws.Activate()
$ws.Range("A45").Select()
